$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B366").Value2 = 'SingleUseId394'
$ws.Range("C366").Value2 = 'Large'
$ws.Range("D366").Value2 = 'Center'
$ws.Range("B367").Value2 = 'SingleUseId396'
$ws.Range("B368").Value2 = 'SingleUseId397'
$ws.Range("D368").Value2 = 'Left'
$ws.Range("E368").Value2 = 'N'
$ws.Range("B369").Value2 = 'SingleUseId398'
$ws.Range("C369").Value2 = 'Bold'
$ws.Range("E369").Value2 = 'ANTISTALL'
$ws.Range("B370").Value2 = 'SingleUseId399'
$ws.Range("C370").Value2 = 'Bold'
$ws.Range("E370").Value2 = 'ANTISTALL'
$ws.Range("B371").Value2 = 'SingleUseId400'
$ws.Range("B372").Value2 = 'SingleUseId401'
$ws.Range("B373").Value2 = 'SingleUseId402'
$ws.Range("C373").Value2 = 'LittleMedium'
$ws.Range("D373").Value2 = 'Center'
$ws.Range("E373").Value2 = 'TRMC 1:'
$ws.Range("B374").Value2 = 'SingleUseId407'
$ws.Range("C374").Value2 = 'Medium'
$ws.Range("D374").Value2 = 'Center'
$ws.Range("E374").Value2 = '<value>'
$ws.Range("B375").Value2 = 'SingleUseId408'
$ws.Range("C375").Value2 = 'Medium'
$ws.Range("D375").Value2 = 'Left'
$c = $ws.Range("E375")
$c.NumberFormat = "@"
$c.Value2 = '0.00'
$c.Style = "Normal"
$ws.Range("F375").Value2 = 'LTR'
$ws.Range("B376").Value2 = 'SingleUseId410'
$ws.Range("C376").Value2 = 'Medium'
$ws.Range("D376").Value2 = 'Center'
$ws.Range("E376").Value2 = '<value>'
$ws.Range("F376").Value2 = 'LTR'
$ws.Range("B377").Value2 = 'SingleUseId411'
$ws.Range("C377").Value2 = 'Medium'
$ws.Range("D377").Value2 = 'Left'
$c = $ws.Range("E377")
$c.NumberFormat = "@"
$c.Value2 = '0.00'
$c.Style = "Normal"
$ws.Range("F377").Value2 = 'LTR'
$ws.Range("B378").Value2 = 'SingleUseId412'
$ws.Range("C378").Value2 = 'Little'
$ws.Range("D378").Value2 = 'Left'
$ws.Range("E378").Value2 = 'VBAT'
$ws.Range("F378").Value2 = 'LTR'
$ws.Range("B379").Value2 = 'SingleUseId414'
$ws.Range("C379").Value2 = 'Little'
$ws.Range("D379").Value2 = 'Left'
$ws.Range("E379").Value2 = 'POIL'
$ws.Range("F379").Value2 = 'LTR'
$ws.Range("B380").Value2 = 'SingleUseId403'
$ws.Range("C380").Value2 = 'Medium'
$ws.Range("D380").Value2 = 'Right'
$ws.Range("E380").Value2 = '<value>'
$ws.Range("F380").Value2 = 'LTR'
$ws.Range("B381").Value2 = 'SingleUseId404'
$ws.Range("C381").Value2 = 'Medium'
$ws.Range("D381").Value2 = 'Left'
$c = $ws.Range("E381")
$c.NumberFormat = "@"
$c.Value2 = '0000'
$c.Style = "Normal"
$ws.Range("F381").Value2 = 'LTR'
